# Daily attendance processing - 2025-12-08 11:49:44
# For every row in column G ("Recorded By"), if the comma-separated list of
# recorders begins with "System", move that "System" entry to the end of
# the list instead of the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newVal = ($rest + "System") -join ", "
            $cell.Value2 = $newVal
        }
    }
}
